$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = [char]39

$ws.Range('D2').Value = '71.073.02'
$ws.Range('E2').Value = '  -2.74%  '
$ws.Range('D3').Value = '3.877.92'
$ws.Range('E3').Value = '  -2.81%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = "$q" + '595.19'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('D6').Value = "$q" + '168.36'
$ws.Range('E6').Value = '  +5.78%  '
$ws.Range('D7').Value = "$q" + '0.673'
$ws.Range('E7').Value = '  -1.81%  '
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('E9').Value = '  +0.40%  '
$ws.Range('D10').Value = "$q" + '0.176'
$ws.Range('E10').Value = '  +3.70%  '
$ws.Range('D11').Value = "$q" + '53.73'
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('D12').Value = "$q" + '0.0000322'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = "$q" + '11.51'
$ws.Range('E13').Value = '  +5.66%  '
$ws.Range('D14').Value = '4.508.16'
$ws.Range('E14').Value = '  -2.48%  '
$ws.Range('D15').Value = '3.883.90'
$ws.Range('E15').Value = '  -2.36%  '
$ws.Range('E16').Value = '  +2.62%  '
$ws.Range('D17').Value = "$q" + '13.90'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('E18').Value = '  -5.59%  '
$ws.Range('E19').Value = '  -2.08%  '
$ws.Range('D20').Value = '71.075.95'
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('D21').Value = "$q" + '437.87'
$ws.Range('E21').Value = '  +0.66%  '
$ws.Range('D22').Value = "$q" + '4.72'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('D23').Value = "$q" + '94.80'
$ws.Range('E23').Value = '  -1.81%  '
$ws.Range('D24').Value = "$q" + '3.28'
$ws.Range('E24').Value = '  -4.81%  '
$ws.Range('D25').Value = "$q" + '13.87'
$ws.Range('E25').Value = '  -3.88%  '
$ws.Range('E26').Value = '  -6.67%  '
$ws.Range('D27').Value = "$q" + '11.36'
$ws.Range('E27').Value = '  +1.57%  '
$ws.Range('D28').Value = "$q" + '5.93'
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('D29').Value = "$q" + '10.36'
$ws.Range('E29').Value = '  -3.37%  '
$ws.Range('D30').Value = "$q" + '35.20'
$ws.Range('E30').Value = '  -3.82%  '
$ws.Range('D31').Value = "$q" + '8.29'
$ws.Range('E31').Value = '  +5.90%  '
$ws.Range('D32').Value = "$q" + '13.63'
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('D33').Value = "$q" + '49.03'
$ws.Range('E33').Value = '  +0.45%  '
$ws.Range('E34').Value = '  -4.28%  '
$ws.Range('D35').Value = "$q" + '70.04'
$ws.Range('E35').Value = '  +1.21%  '
$ws.Range('D36').Value = '0.0₃0994'
$ws.Range('E36').Value = '  +12.94%  '
$ws.Range('D37').Value = "$q" + '634.43'
$ws.Range('E37').Value = '  -6.57%  '
$ws.Range('D38').Value = "$q" + '0.428'
$ws.Range('E38').Value = '  -2.29%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = "$q" + '0.146'
$ws.Range('E39').Value = '  -0.64%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').Value = "$q" + '1.00'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = "$q" + '3.34'
$ws.Range('E41').Value = '  +27.24%  '
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('D43').Value = "$q" + '3.29'
$ws.Range('E43').Value = '  -2.01%  '
$ws.Range('D44').Value = "$q" + '0.0472'
$ws.Range('E44').Value = '  -3.36%  '
$ws.Range('D45').Value = "$q" + '10.19'
$ws.Range('E45').Value = '  -6.05%  '
$ws.Range('D46').Value = "$q" + '2.72'
$ws.Range('E46').Value = '  +1.60%  '
$ws.Range('E47').Value = '  -3.85%  '
$ws.Range('E48').Value = '  -15.51%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.869.52'
$ws.Range('E49').Value = '  +2.62%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = "$q" + '3.30'
$ws.Range('E50').Value = '  -3.49%  '
$ws.Range('D51').Value = "$q" + '0.000273'
$ws.Range('E51').Value = '  +0.68%  '
